$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p145r_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p145r_2</id>", 2)
$d.Content.Find.Execute("<id>p145v_1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p145v_1</id>", 2)
$d.Content.Find.Execute("<id>p145v_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p145v_2</id>", 2)
